# The sheet "Current-Dollar GDP" lists US states with their GDP figure in
# column B. Row 2 held the "United States" total (bold, out of alphabetical
# order); every other row (3-53) was already sorted alphabetically by state
# name. This script performs the equivalent of re-sorting A2:B53 ascending
# by column A, which moves "United States" down to sit between "Texas" and
# "Utah" (row 46), shifting Alabama..Texas up by one row (rows 2-45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the cell formatting between row 2 and row 46 ------------------
# Row 2 currently carries the bold "total" style; row 46 carries the regular
# "state" style. After the sort those two rows trade places, so swap their
# formats (copy/paste-format keeps the original style indices instead of
# synthesising new ones).
$ws.Range("A46:B46").Copy() | Out-Null
$ws.Range("A2:B2").Copy() | Out-Null

$tempRow = 60
$ws.Range("A2:B2").Copy()
$ws.Range("A" + $tempRow + ":B" + $tempRow).PasteSpecial(-4122) | Out-Null   # xlPasteFormats, stash row2's bold format

$ws.Range("A46:B46").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122) | Out-Null                            # row2 <- row46's regular format

$ws.Range("A" + $tempRow + ":B" + $tempRow).Copy()
$ws.Range("A46:B46").PasteSpecial(-4122) | Out-Null                          # row46 <- stashed bold format

$ws.Range("A" + $tempRow + ":B" + $tempRow).Clear() | Out-Null
$excel.CutCopyMode = $false

# --- 2. Write the re-sorted values -----------------------------------------
$rows = @(
    @{Row=2; Label='Alabama '; Value=248543},
    @{Row=3; Label='Alaska '; Value=55531},
    @{Row=4; Label='Arizona '; Value=409577},
    @{Row=5; Label='Arkansas'; Value=146292},
    @{Row=6; Label='California '; Value=3353473},
    @{Row=7; Label='Colorado '; Value=425595},
    @{Row=8; Label='Connecticut'; Value=299819},
    @{Row=9; Label='Delaware '; Value=81019},
    @{Row=10; Label='District of Columbia '; Value=153979},
    @{Row=11; Label='Florida '; Value=1226956},
    @{Row=12; Label='Georgia '; Value=687680},
    @{Row=13; Label='Hawaii '; Value=92445},
    @{Row=14; Label='Idaho '; Value=93785},
    @{Row=15; Label='Illinois '; Value=956366},
    @{Row=16; Label='Indiana'; Value=423052},
    @{Row=17; Label='Iowa '; Value=225144},
    @{Row=18; Label='Kansas'; Value=196495},
    @{Row=19; Label='Kentucky'; Value=239208},
    @{Row=20; Label='Louisiana '; Value=256905},
    @{Row=21; Label='Maine '; Value=76045},
    @{Row=22; Label='Maryland '; Value=451635},
    @{Row=23; Label='Massachusetts '; Value=637424},
    @{Row=24; Label='Michigan '; Value=567450},
    @{Row=25; Label='Minnesota '; Value=415393},
    @{Row=26; Label='Mississippi '; Value=126155},
    @{Row=27; Label='Missouri '; Value=365899},
    @{Row=28; Label='Montana '; Value=59129},
    @{Row=29; Label='Nebraska'; Value=153586},
    @{Row=30; Label='Nevada'; Value=192228},
    @{Row=31; Label='New Hampshire '; Value=94305},
    @{Row=32; Label='New Jersey'; Value=686842},
    @{Row=33; Label='New Mexico'; Value=110696},
    @{Row=34; Label='New York'; Value=1893670},
    @{Row=35; Label='North Carolina '; Value=659616},
    @{Row=36; Label='North Dakota '; Value=66632},
    @{Row=37; Label='Ohio '; Value=745000},
    @{Row=38; Label='Oklahoma'; Value=211053},
    @{Row=39; Label='Oregon '; Value=270120},
    @{Row=40; Label='Pennsylvania '; Value=849818},
    @{Row=41; Label='Rhode Island '; Value=65348},
    @{Row=42; Label='South Carolina '; Value=271374},
    @{Row=43; Label='South Dakota '; Value=61732},
    @{Row=44; Label='Tennessee '; Value=421086},
    @{Row=45; Label='Texas'; Value=2006662},
    @{Row=46; Label='       United States'; Value=23202344}
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Label
    $ws.Cells.Item($r, 2).Value = $entry.Value
}

# --- 3. Restore the selection left behind by the sort dialog ---------------
$ws.Range("A2:B2").Select() | Out-Null
